$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.538.20'
$ws.Range('D2').Style = "Normal"
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.23'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.78%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.79'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.599'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.53'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.154'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +3.13%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.090.24'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.34'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +13.11%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '60.522.62'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.633.74'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.87%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.57'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.02%  '
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '348.55'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.90'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.32%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.530'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.91'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('E27').Value = '  +5.92%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.05'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +11.92%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0798'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.69'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +5.41%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '170.10'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +5.64%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.07'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +11.20%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.44'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.14%  '
$ws.Range('E36').Value = '  +8.36%  '
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '330.64'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +12.41%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.01'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.17%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '38.65'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.18%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.873'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.61%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.18'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +7.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '20.73'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.71%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0999'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.99%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '132.84'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.89%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '20.08'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0557'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.609'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('E50').Value = '  +2.51%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '10.73'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.53%  '
